$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 904.2857
$ws.Range("I46").Value = 498.57144
$ws.Range("J46").Value = 1107.1428
$ws.Range("K46").Value = 1495.71432
$ws.Range("L46").Value = 3321.4284
$ws.Range("M46").Value = -1376.71432
$ws.Range("N46").Value = -3559.4284
$ws.Range("H60").Value = 904.2857
$ws.Range("I60").Value = 498.57144
$ws.Range("J60").Value = 1107.1428
$ws.Range("K60").Value = 1495.71432
$ws.Range("L60").Value = 3321.4284
$ws.Range("M60").Value = -1011.71432
$ws.Range("N60").Value = -4289.428400000001
$ws.Range("H69").Value = 34493640
$ws.Range("J69").Value = 35721984
$ws.Range("L69").Value = 107165952
$ws.Range("N69").Value = -107167700
$ws.Range("H72").Value = 34493640
$ws.Range("J72").Value = 35721984
$ws.Range("L72").Value = 321497856
$ws.Range("N72").Value = -321506592
$ws.Range("H96").Value = 592.6
$ws.Range("I96").Value = 633.44446
$ws.Range("K96").Value = 1900.33338
$ws.Range("M96").Value = -527.33338
$ws.Range("H113").Value = 6137.0713
$ws.Range("I113").Value = 2995
$ws.Range("J113").Value = 6378.769
$ws.Range("K113").Value = 2995
$ws.Range("L113").Value = 6378.769
$ws.Range("M113").Value = 259
$ws.Range("N113").Value = -12886.769
$ws.Range("H115").Value = 332.8
$ws.Range("I115").Value = 332.8
$ws.Range("K115").Value = 998.4000000000001
$ws.Range("M115").Value = 568.5999999999999
$ws.Range("H138").Value = 2465.1929
$ws.Range("J138").Value = 4862.2383
$ws.Range("L138").Value = 14586.7149
$ws.Range("N138").Value = -24866.7149
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 80000
$ws.Range("J24").Value = 80000
$ws.Range("L24").Value = 80000
$ws.Range("N24").Value = -80748
$ws.Range("H43").Value = 28332.666
$ws.Range("J43").Value = 29999
$ws.Range("L43").Value = 29999
$ws.Range("N43").Value = -30625
$ws.Range("H45").Value = 2204397.5
$ws.Range("I45").Value = 3667664.8
$ws.Range("J45").Value = 9497
$ws.Range("K45").Value = 3667664.8
$ws.Range("L45").Value = 9497
$ws.Range("M45").Value = -3667287.8
$ws.Range("N45").Value = -10251
$ws.Range("H96").Value = 51035.43
$ws.Range("J96").Value = 51035.43
$ws.Range("L96").Value = 51035.43
$ws.Range("N96").Value = -56527.43
$ws.Range("H97").Value = 10468.275
$ws.Range("I97").Value = 8307.125
$ws.Range("J97").Value = 20841.8
$ws.Range("K97").Value = 8307.125
$ws.Range("L97").Value = 20841.8
$ws.Range("M97").Value = -7811.125
$ws.Range("N97").Value = -21833.8
$ws.Range("H100").Value = 80000
$ws.Range("J100").Value = 80000
$ws.Range("L100").Value = 80000
$ws.Range("N100").Value = -82164
$ws.Range("H122").Value = 3027.963
$ws.Range("I122").Value = 2832.65
$ws.Range("J122").Value = 3586
$ws.Range("K122").Value = 8497.950000000001
$ws.Range("L122").Value = 10758
$ws.Range("M122").Value = -6047.950000000001
$ws.Range("N122").Value = -15658
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3268.6316
$ws.Range("I99").Value = 2055.3635
$ws.Range("J99").Value = 4936.875
$ws.Range("K99").Value = 2055.3635
$ws.Range("L99").Value = 4936.875
$ws.Range("M99").Value = -557.3634999999999
$ws.Range("N99").Value = -7932.875
$ws.Range("H105").Value = 1711.25
$ws.Range("I105").Value = 1711.25
$ws.Range("K105").Value = 1711.25
$ws.Range("M105").Value = 35.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 32026.967
$ws.Range("I132").Value = 2500.75
$ws.Range("J132").Value = 445394
$ws.Range("K132").Value = 7502.25
$ws.Range("L132").Value = 1336182
$ws.Range("M132").Value = -4972.25
$ws.Range("N132").Value = -1341242
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2563.6128
$ws.Range("I113").Value = 3298.0908
$ws.Range("J113").Value = 2159.65
$ws.Range("K113").Value = 9894.2724
$ws.Range("L113").Value = 6478.950000000001
$ws.Range("M113").Value = -7724.2724
$ws.Range("N113").Value = -10818.95
$ws.Range("H122").Value = 880.9167
$ws.Range("I122").Value = 722.25
$ws.Range("K122").Value = 6500.25
$ws.Range("M122").Value = -4050.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2530.0908
$ws.Range("I102").Value = 2019.96
$ws.Range("K102").Value = 2019.96
$ws.Range("M102").Value = -397.96
$ws.Range("H113").Value = 1896.6666
$ws.Range("I113").Value = 1896.6666
$ws.Range("K113").Value = 1896.6666
$ws.Range("M113").Value = 273.3334
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3946.8235
$ws.Range("I7").Value = 1879
$ws.Range("J7").Value = 8909.6
$ws.Range("K7").Value = 1879
$ws.Range("L7").Value = 8909.6
$ws.Range("M7").Value = -1767
$ws.Range("N7").Value = -9133.6
$ws.Range("H40").Value = 7170.357
$ws.Range("I40").Value = 5820.778
$ws.Range("J40").Value = 9599.6
$ws.Range("K40").Value = 5820.778
$ws.Range("L40").Value = 9599.6
$ws.Range("M40").Value = -5684.778
$ws.Range("N40").Value = -9871.6
$ws.Range("H109").Value = 59987.668
$ws.Range("J109").Value = 59987.668
$ws.Range("L109").Value = 59987.668
$ws.Range("N109").Value = -62761.668
$ws.Range("H122").Value = 6068.3687
$ws.Range("I122").Value = 4652.1
$ws.Range("K122").Value = 13956.3
$ws.Range("M122").Value = -11506.3
$ws.Range("H126").Value = 3946.8235
$ws.Range("I126").Value = 1879
$ws.Range("J126").Value = 8909.6
$ws.Range("K126").Value = 5637
$ws.Range("L126").Value = 26728.8
$ws.Range("M126").Value = -3167
$ws.Range("N126").Value = -31668.8
$ws.Range("H132").Value = 8383.817999999999
$ws.Range("I132").Value = 13401
$ws.Range("J132").Value = 4910.385
$ws.Range("K132").Value = 40203
$ws.Range("L132").Value = 14731.155
$ws.Range("M132").Value = -37673
$ws.Range("N132").Value = -19791.155
$ws.Range("H136").Value = 40646.594
$ws.Range("I136").Value = 61385.883
$ws.Range("K136").Value = 184157.649
$ws.Range("M136").Value = -181607.649
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2773.6155
$ws.Range("I126").Value = 2644.56
$ws.Range("K126").Value = 7933.68
$ws.Range("M126").Value = -5463.68
